# Selenium_Java_Allure/src/test/resources/data.xlsx update
# Expands the Data Name / Data Value table from 10 data rows to 28 data
# rows of test data, bolds the header row, and tweaks a couple of value
# cells' formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Full Data Name / Data Value table (row 1 = header) -------------------
$data = @(
  @("Data Name", "Data Value "),
  @("chargeItemValue1", "UCCITest#01"),
  @("chargeItemValue2", "UCCITest#02"),
  @("chargeItemValue3", "UCCITest#03"),
  @("chargeItemValue4", "UC_CITest#04"),
  @("chargeItemValue5", "UCCITest#05"),
  @("chargeItemValue6", "UCCITest#06"),
  @("chargeItemValue7", "UCCITest#07"),
  @("chargeItemValue8", "UCCITest#08"),
  @("chargeItemValue9", "UCCITest#09"),
  @("chargeItemValue10", "UCCITest#10"),
  @("chargeItemValue11", "UCCITest#11"),
  @("chargeItemValue12", "UCCITest#12"),
  @("chargeItemValue13", "UCCITest#13"),
  @("chargeItemDescription1", "For testing"),
  @("chargeItemDescription2", "FortestingFortestingFortestingFortestingFortestingFortestingFortestingFortestingFortestingFortestingFortestingFortestingFortestingFortestingFortestingFortestingFortestingFortestingFortestingFortestingF"),
  @("amount1", '"500"'),
  @("amount2", '"600.50"'),
  @("searchItemValue1", "UCCITest#01"),
  @("searchItemValue2", "UC_CI_Test#02"),
  @("expectedValue1", "UC_CI_Test#02"),
  @("expectedErrorValue1", "There is an existing charge item with the same name."),
  @("expectedErrorValue2", "You are about to cancel your entry and your configurations will not be saved. Are you sure you want to proceed?"),
  @("expectedErrorValue3", "69"),
  @("expectedLength1", "200/200"),
  @("mandatoryErrorValue1", "Enter a value to proceed."),
  @("mandatoryErrorValue2", "Select a value to proceed."),
  @("existingChargeItemValue1", "UC_CI_Test#01"),
  @("createDrawerTitleValue", "Create charge item")
)

for ($i = 0; $i -lt $data.Length; $i++) {
  $r = $i + 1
  $ws.Cells.Item($r, 1).Value = $data[$i][0]
  $ws.Cells.Item($r, 2).Value = $data[$i][1]
}

# Header row -> bold
$ws.Range("A1:B1").Font.Bold = $true

# The very long description value keeps its default (left/general) alignment
# explicitly re-applied, which is how it ends up carrying its own style
# record in the workbook.
$ws.Range("B16").WrapText = $false

# The numeric-looking "69" value is forced to text with a leading
# apostrophe so Excel keeps it left-aligned as a quoted string.
$ws.Range("B24").Value = "'69"

# Column sizing: the Data Name column is re-fit to its (now longer)
# content, while the Data Value column is widened by hand so the long
# test strings stay readable.
$ws.Columns("A").ColumnWidth = 23.35
$ws.Columns("B").ColumnWidth = 22.65

# Leave selection where the author ended up after entering the data.
$ws.Range("H31").Select()
